$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.004791757726395753
$ws.Range("C2").Value = 0.0007813039001524819
$ws.Range("B3").Value = 0.000608646277943711
$ws.Range("C3").Value = 0.0002040837046817845
$ws.Range("A4").Value = "GOOG"
$ws.Range("B4").Value = 0.0001578400726809109
$ws.Range("C4").Value = 0.000157536931257989
$ws.Range("A5").Value = "GOOGL"
$ws.Range("B5").Value = 0.0001695143254891404
$ws.Range("C5").Value = 0.0001594834729978914
$ws.Range("B7").Value = 0.0004710739231839039
$ws.Range("C7").Value = 0.0003049396500722201
$ws.Range("B8").Value = 0.0009251793916008272
$ws.Range("C8").Value = 0.0003493147409145823
$ws.Range("B10").Value = 0.00005969832236319824
$ws.Range("C10").Value = 0.0001216729727845804
$ws.Range("B11").Value = 0.001046814651255135
$ws.Range("C11").Value = 0.0002793057785124484
$ws.Range("B14").Value = 0.001552509833639477
$ws.Range("C14").Value = 0.0003080525190485758
$ws.Range("A15").Value = "PEP"
$ws.Range("B15").Value = 0.0003452022700634803
$ws.Range("C15").Value = 0.000071983828376838
$ws.Range("A16").Value = "AZN"
$ws.Range("B16").Value = -0.0006208107451509798
$ws.Range("C16").Value = 0.0002472116218947619
$ws.Range("B17").Value = 0.0006712548222779616
$ws.Range("C17").Value = 0.0001385912116243684
$ws.Range("B19").Value = 0.0006420875880178026
$ws.Range("C19").Value = 0.0001811539506805654
$ws.Range("B20").Value = 0.001187964627932078
$ws.Range("C20").Value = 0.0003132064104617103
$ws.Range("B21").Value = 0.001280131901946855
$ws.Range("C21").Value = 0.0001753747607148026
$ws.Range("A22").Value = "AMGN"
$ws.Range("B22").Value = -0.0002138669875757005
$ws.Range("C22").Value = 0.0002375851329806968
$ws.Range("A23").Value = "INTU"
$ws.Range("B23").Value = 0.0008204261379338005
$ws.Range("C23").Value = 0.0001541026309202054
$ws.Range("B25").Value = 0.002288414744028531
$ws.Range("C25").Value = 0.0003758874907965164
$ws.Range("B26").Value = 0.0009270755367825428
$ws.Range("C26").Value = 0.0001144594514881142
$ws.Range("B27").Value = 0.0006499838885816938
$ws.Range("C27").Value = 0.0003507276499620088
$ws.Range("B28").Value = 0.0005934713532669781
$ws.Range("C28").Value = 0.000125638818668137
$ws.Range("A31").Value = "MU"
$ws.Range("B31").Value = 0.001693402989727995
$ws.Range("C31").Value = 0.001029879002263344
$ws.Range("A32").Value = "ADP"
$ws.Range("B32").Value = 0.000992937478419769
$ws.Range("C32").Value = 0.0001108287177409622
$ws.Range("B33").Value = 0.001258298230554382
$ws.Range("C33").Value = 0.0002204926327720668
$ws.Range("B35").Value = -0.0001327435953880019
$ws.Range("C35").Value = 0.0001572055422603591
$ws.Range("A36").Value = "LRCX"
$ws.Range("B36").Value = 0.001554779270507288
$ws.Range("C36").Value = 0.0003175188223807012
$ws.Range("A37").Value = "KLAC"
$ws.Range("B37").Value = 0.0006617849538736269
$ws.Range("C37").Value = 0.000176377742102125
$ws.Range("A38").Value = "MELI"
$ws.Range("B38").Value = 0.001414183280644268
$ws.Range("C38").Value = 0.0005015673650968572
$ws.Range("A39").Value = "GILD"
$ws.Range("B39").Value = -0.001159193092263775
$ws.Range("C39").Value = 0.0002872735440944644
$ws.Range("B40").Value = 0.0003879786706433413
$ws.Range("C40").Value = 0.0002042997220563958
$ws.Range("B41").Value = 0.0001168621770919631
$ws.Range("C41").Value = 0.0002846239921258829
$ws.Range("B42").Value = 0.001280820366687428
$ws.Range("C42").Value = 0.0001526327418421319
$ws.Range("A43").Value = "SNPS"
$ws.Range("B43").Value = 0.001076807426562563
$ws.Range("C43").Value = 0.0001280042808616382
$ws.Range("A44").Value = "PYPL"
$ws.Range("B44").Value = 0.0005074157314395192
$ws.Range("C44").Value = 0.0003085645244898223
$ws.Range("B46").Value = 0.001000704087030286
$ws.Range("C46").Value = 0.0002465631699095578
$ws.Range("B48").Value = 0.001424263413040009
$ws.Range("C48").Value = 0.0002638560309284628
$ws.Range("A49").Value = "WDAY"
$ws.Range("B49").Value = -0.0006478477830791533
$ws.Range("C49").Value = 0.0007386675435605258
$ws.Range("A50").Value = "MRVL"
$ws.Range("B50").Value = 0.001894871010761541
$ws.Range("C50").Value = 0.0005145473683505556
$ws.Range("A51").Value = "FTNT"
$ws.Range("B51").Value = -0.00001452342082719582
$ws.Range("C51").Value = 0.0004721799121230202
$ws.Range("A52").Value = "NXPI"
$ws.Range("B52").Value = 0.0005937381103851151
$ws.Range("C52").Value = 0.0005494577051172493
$ws.Range("B54").Value = -0.00004713989730441069
$ws.Range("C54").Value = 0.0001640338333908919
$ws.Range("B55").Value = 0.001690814103568938
$ws.Range("C55").Value = 0.0006042643448072175
$ws.Range("B56").Value = 0.001314377003276264
$ws.Range("C56").Value = 0.0002610052993351169
$ws.Range("A57").Value = "CPRT"
$ws.Range("B57").Value = 0.001537434689359692
$ws.Range("C57").Value = 0.000184575384017428
$ws.Range("A59").Value = "AEP"
$ws.Range("B59").Value = 0.0004447342210871816
$ws.Range("C59").Value = 0.0001255333303717242
$ws.Range("A60").Value = "MNST"
$ws.Range("B60").Value = -0.0003253870869810081
$ws.Range("C60").Value = 0.0002874577754103618
$ws.Range("A61").Value = "PAYX"
$ws.Range("B61").Value = 0.0007944879934807374
$ws.Range("C61").Value = 0.0001267022309677935
$ws.Range("A62").Value = "KDP"
$ws.Range("B62").Value = 0.00002486416858565427
$ws.Range("C62").Value = 0.0001231958349491246
$ws.Range("A63").Value = "TEAM"
$ws.Range("B63").Value = -0.0005160409859583973
$ws.Range("C63").Value = 0.0008129416636759335
$ws.Range("A64").Value = "ROST"
$ws.Range("B64").Value = 0.0007842419758300012
$ws.Range("C64").Value = 0.0001762896607575369
$ws.Range("A65").Value = "FAST"
$ws.Range("B65").Value = 0.0007393852622347126
$ws.Range("C65").Value = 0.000237101092428346
$ws.Range("A66").Value = "ODFL"
$ws.Range("B66").Value = 0.001530852491970737
$ws.Range("C66").Value = 0.000284365347228002
$ws.Range("A67").Value = "KHC"
$ws.Range("B67").Value = 0.0008422548716767359
$ws.Range("C67").Value = 0.0001654461098804011
$ws.Range("A68").Value = "MCHP"
$ws.Range("B68").Value = 0.001435790485856591
$ws.Range("C68").Value = 0.0002541671089951331
$ws.Range("B69").Value = 0.001099212824391787
$ws.Range("C69").Value = 0.0001990675577491658
$ws.Range("B71").Value = 0.0002990666476850587
$ws.Range("C71").Value = 0.0001148203702621627
$ws.Range("B72").Value = 0.0006981535920071758
$ws.Range("C72").Value = 0.0003375210087248839
$ws.Range("B73").Value = -0.0001465932538590275
$ws.Range("C73").Value = 0.0003799404094710131
$ws.Range("B74").Value = 0.001325896358313684
$ws.Range("C74").Value = 0.0006030836269188352
$ws.Range("B75").Value = -0.0003398158017141503
$ws.Range("C75").Value = 0.0001958379840393211
$ws.Range("A76").Value = "LULU"
$ws.Range("B76").Value = 0.0006031267234766266
$ws.Range("C76").Value = 0.0004333103936674401
$ws.Range("A77").Value = "XEL"
$ws.Range("B77").Value = 0.0006535408186765775
$ws.Range("C77").Value = 0.000123123541735914
$ws.Range("B80").Value = 0.001001691522343329
$ws.Range("C80").Value = 0.0002214376871376766
$ws.Range("A83").Value = "SMCI"
$ws.Range("B83").Value = 0.0006797830151301731
$ws.Range("C83").Value = 0.001316613247760909
$ws.Range("A85").Value = "TTWO"
$ws.Range("B85").Value = 0.001459654289844646
$ws.Range("C85").Value = 0.0003452783591290246
